$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns D and E: copy C1's format (bold/border/centered header
# style) onto D1:E1, then set their values.
$ws.Range("C1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4

# Updated predicted values for columns C, D, E (rows 2-10)
$ws.Range("C2").Value = -5.030479892299043
$ws.Range("D2").Value = -4.700423608398296
$ws.Range("E2").Value = -4.336375677541507

$ws.Range("C3").Value = -1.153235334261761
$ws.Range("D3").Value = -1.154630231040654
$ws.Range("E3").Value = -1.134646576557121

$ws.Range("C4").Value = -0.07292569980107731
$ws.Range("D4").Value = 0.02093359320803632
$ws.Range("E4").Value = 0.09358456433631321

$ws.Range("C5").Value = -0.417749988516372
$ws.Range("D5").Value = -0.1565149916917322
$ws.Range("E5").Value = 0.05089419996022063

$ws.Range("C6").Value = 0.01457436480836208
$ws.Range("D6").Value = -0.04514720434098682
$ws.Range("E6").Value = -0.09996768479672546

$ws.Range("C7").Value = 0.1039564587721915
$ws.Range("D7").Value = 0.06550401096984124
$ws.Range("E7").Value = 0.02827358676505631

$ws.Range("C8").Value = 0.1352692197136115
$ws.Range("D8").Value = 0.0397665396836725
$ws.Range("E8").Value = -0.05068768714866342

$ws.Range("C9").Value = 0.02838054686518928
$ws.Range("D9").Value = 0.024635519709972
$ws.Range("E9").Value = 0.02148294890355889

$ws.Range("C10").Value = 0.02323175602225529
$ws.Range("D10").Value = 0.01754457296797431
$ws.Range("E10").Value = 0.01323784340406731
